{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst newValues = [\n  \"22+37=59\",\n  \"69+30=99\",\n  \"3+48=51\",\n  \"62-17=45\",\n  \"37-1=36\",\n  \"77-43=34\",\n  \"96-13=83\",\n  \"75-48=27\",\n  \"74-20=54\",\n  \"92-78=14\",\n  \"42-5=37\",\n  \"51-24=27\",\n  \"37+40=77\",\n  \"94-29=65\",\n  \"25+0=25\",\n  \"21-18=3\",\n  \"73-5=68\",\n  \"6+66=72\",\n  \"94-22=72\",\n  \"52+31=83\",\n  \"84+0=84\",\n  \"44+10=54\",\n  \"5+14=19\",\n  \"70-56=14\",\n  \"47-30=17\",\n  \"84-42=42\",\n  \"51+7=58\",\n  \"46+44=90\",\n  \"15+34=49\",\n  \"87-45=42\",\n  \"31-0=31\",\n  \"4+38=42\",\n  \"68-22=46\",\n  \"88+1=89\",\n  \"23+16=39\",\n  \"33-22=11\",\n  \"22+31=53\",\n  \"42+57=99\",\n  \"10+45=55\",\n  \"79-26=53\",\n  \"40-13=27\",\n  \"75-22=53\",\n  \"70+14=84\",\n  \"15+37=52\",\n  \"82-16=66\",\n  \"81-6=75\",\n  \"35+22=57\",\n  \"43-39=4\",\n  \"80-12=68\",\n  \"59-12=47\",\n  \"93-56=37\",\n  \"54-45=9\",\n  \"85-57=28\",\n  \"22+25=47\",\n  \"86-59=27\",\n  \"42+42=84\",\n  \"91-11=80\",\n  \"43-36=7\",\n  \"13+3=16\",\n  \"19+80=99\",\n  \"13-6=7\",\n  \"6+8=14\",\n  \"38+20=58\",\n  \"63+33=96\",\n  \"69-30=39\",\n  \"85+13=98\",\n  \"21-17=4\",\n  \"21+25=46\",\n  \"13+25=38\",\n  \"7+44=51\",\n  \"66+30=96\",\n  \"21+70=91\",\n  \"4+7=11\",\n  \"59+35=94\",\n  \"39+55=94\",\n  \"64-25=39\",\n  \"86-71=15\",\n  \"69+1=70\",\n  \"34-10=24\",\n  \"48-47=1\",\n  \"95-57=38\",\n  \"41-23=18\",\n  \"82-1=81\",\n  \"1+45=46\",\n  \"71-70=1\",\n  \"49-31=18\",\n  \"1+32=33\",\n  \"67-28=39\",\n  \"57-5=52\",\n  \"39-22=17\",\n  \"23-6=17\",\n  \"60-16=44\",\n  \"59-18=41\",\n  \"89-7=82\",\n  \"11-10=1\",\n  \"72-59=13\",\n  \"95-84=11\",\n  \"10+21=31\",\n  \"34+0=34\",\n  \"0+16=16\"\n];\nconst columnCount = 5;\nlet idx = 0;\nfor (let row = 0; row < table.rowCount; row++) {\n  for (let col = 0; col < columnCount; col++) {\n    const cell = table.getCell(row, col);\n    cell.value = newValues[idx];\n    idx++;\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$newValues = @(\n    \"22+37=59\",\n    \"69+30=99\",\n    \"3+48=51\",\n    \"62-17=45\",\n    \"37-1=36\",\n    \"77-43=34\",\n    \"96-13=83\",\n    \"75-48=27\",\n    \"74-20=54\",\n    \"92-78=14\",\n    \"42-5=37\",\n    \"51-24=27\",\n    \"37+40=77\",\n    \"94-29=65\",\n    \"25+0=25\",\n    \"21-18=3\",\n    \"73-5=68\",\n    \"6+66=72\",\n    \"94-22=72\",\n    \"52+31=83\",\n    \"84+0=84\",\n    \"44+10=54\",\n    \"5+14=19\",\n    \"70-56=14\",\n    \"47-30=17\",\n    \"84-42=42\",\n    \"51+7=58\",\n    \"46+44=90\",\n    \"15+34=49\",\n    \"87-45=42\",\n    \"31-0=31\",\n    \"4+38=42\",\n    \"68-22=46\",\n    \"88+1=89\",\n    \"23+16=39\",\n    \"33-22=11\",\n    \"22+31=53\",\n    \"42+57=99\",\n    \"10+45=55\",\n    \"79-26=53\",\n    \"40-13=27\",\n    \"75-22=53\",\n    \"70+14=84\",\n    \"15+37=52\",\n    \"82-16=66\",\n    \"81-6=75\",\n    \"35+22=57\",\n    \"43-39=4\",\n    \"80-12=68\",\n    \"59-12=47\",\n    \"93-56=37\",\n    \"54-45=9\",\n    \"85-57=28\",\n    \"22+25=47\",\n    \"86-59=27\",\n    \"42+42=84\",\n    \"91-11=80\",\n    \"43-36=7\",\n    \"13+3=16\",\n    \"19+80=99\",\n    \"13-6=7\",\n    \"6+8=14\",\n    \"38+20=58\",\n    \"63+33=96\",\n    \"69-30=39\",\n    \"85+13=98\",\n    \"21-17=4\",\n    \"21+25=46\",\n    \"13+25=38\",\n    \"7+44=51\",\n    \"66+30=96\",\n    \"21+70=91\",\n    \"4+7=11\",\n    \"59+35=94\",\n    \"39+55=94\",\n    \"64-25=39\",\n    \"86-71=15\",\n    \"69+1=70\",\n    \"34-10=24\",\n    \"48-47=1\",\n    \"95-57=38\",\n    \"41-23=18\",\n    \"82-1=81\",\n    \"1+45=46\",\n    \"71-70=1\",\n    \"49-31=18\",\n    \"1+32=33\",\n    \"67-28=39\",\n    \"57-5=52\",\n    \"39-22=17\",\n    \"23-6=17\",\n    \"60-16=44\",\n    \"59-18=41\",\n    \"89-7=82\",\n    \"11-10=1\",\n    \"72-59=13\",\n    \"95-84=11\",\n    \"10+21=31\",\n    \"34+0=34\",\n    \"0+16=16\"\n)\n$columnCount = 5\n$idx = 0\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    for ($c = 1; $c -le $columnCount; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $newValues[$idx]\n        $idx++\n    }\n}"}
